$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P1").Value = 'construct'
$ws.Range("P2").Value = 'Ctcf_CDS_1'
$ws.Range("P3").Value = 'Rnf2_CDS_1'
$ws.Range("P4").Value = 'Ehmt2_CDS_1'
$ws.Range("P5").Value = 'Parp1_CDS_1'
$ws.Range("P6").Value = 'Parp1_CDS_2'
$ws.Range("P7").Value = 'Zfp57_CDS_1'
$ws.Range("P8").Value = 'Atf7ip_CDS_1'
$ws.Range("P9").Value = 'Hdac2_CDS_1'
$ws.Range("P10").Value = 'Kmt5a_CDS_1'
$ws.Range("P11").Value = 'Suv39h1_CDS_1'
$ws.Range("P12").Value = 'Smchd1_CDS_1'
$ws.Range("P13").Value = 'Sirt1_CDS_1'
$ws.Range("P14").Value = 'Atrx_CDS_1'
$ws.Range("P15").Value = 'Rad21_CDS_1'
$ws.Range("P16").Value = 'Eed_3UTR_1'
$ws.Range("P17").Value = 'Hdac1_CDS_1'
$ws.Range("P18").Value = 'Dnmt1_CDS_1'
$ws.Range("P19").Value = 'Eed_CDS_2'
$ws.Range("P20").Value = 'Rrm2_CDS_1'
$ws.Range("P21").Value = 'Atrx_CDS_2'
$ws.Range("P22").Value = 'Dnmt1_CDS_2'
$ws.Range("P23").Value = 'Sirt1_CDS_2'
$ws.Range("P24").Value = 'Parp1_3UTR_3'
$ws.Range("P25").Value = 'Hdac2_CDS_2'
$ws.Range("P26").Value = 'Smchd1_CDS_2'
$ws.Range("P27").Value = 'Ehmt2_CDS_2'
$ws.Range("P28").Value = 'Rrm2_CDS_2'
$ws.Range("P29").Value = 'Rrm2_CDS_3'
$ws.Range("P30").Value = 'Smchd1_CDS_3'
$ws.Range("P31").Value = 'Rrm2_CDS_4'
$ws.Range("P32").Value = 'Kmt5a_CDS_2'
$ws.Range("P33").Value = 'LUCIFERASE_CDS'
$ws.Range("P34").Value = 'Dnmt1_CDS_3'
$ws.Range("P35").Value = 'Dnmt1_3UTR_4'
$ws.Range("P36").Value = 'Rad21_CDS_2'
$ws.Range("P37").Value = 'Ring1_CDS_1'
$ws.Range("P38").Value = 'Ehmt2_3UTR_3'
$ws.Range("P39").Value = 'Rrm2_CDS_5'
$ws.Range("P40").Value = 'Ring1_CDS_2'
$ws.Range("P41").Value = 'Kmt5a_CDS_3'
$ws.Range("P42").Value = 'Ehmt2_CDS_4'
$ws.Range("P43").Value = 'Ctcf_CDS_2'
$ws.Range("P44").Value = 'Suv39h1_CDS_2'
$ws.Range("P45").Value = 'Atrx_CDS_3'
$ws.Range("P46").Value = 'GFP_CDS'
$ws.Range("P47").Value = 'Hdac2_CDS_3'
$ws.Range("P48").Value = 'Hdac2_CDS_4'
$ws.Range("P49").Value = 'Smchd1_CDS_4'
$ws.Range("P50").Value = 'Eed_CDS_3'
$ws.Range("P51").Value = 'Atf7ip_CDS_2'
$ws.Range("P52").Value = 'Eed_CDS_4'
$ws.Range("P53").Value = 'Rad21_3UTR_3'
$ws.Range("P54").Value = 'Hdac1_3UTR_2'
$ws.Range("P55").Value = 'Rad21_CDS_4'
$ws.Range("P56").Value = 'Ehmt2_CDS_5'
$ws.Range("P57").Value = 'Rnf2_CDS_2'
$ws.Range("P58").Value = 'Kmt5a_CDS_4'
$ws.Range("P59").Value = 'Sirt1_CDS_3'
$ws.Range("P60").Value = 'Zfp57_CDS_2'
$ws.Range("P61").Value = 'Ring1_CDS_3'
$ws.Range("P62").Value = 'Hdac1_CDS_3'
$ws.Range("P63").Value = 'Sirt1_CDS_4'
$ws.Range("P64").Value = 'Atf7ip_CDS_3'
$ws.Range("P65").Value = 'Parp1_CDS_4'
$ws.Range("P66").Value = 'Dnmt1_CDS_5'
$ws.Range("P67").Value = 'Ctcf_CDS_3'
$ws.Range("P68").Value = 'Zfp57_CDS_3'
$ws.Range("P69").Value = 'Ctcf_CDS_4'
$ws.Range("P70").Value = 'Rad21_CDS_5'
$ws.Range("P71").Value = 'Rnf2_3UTR_3'
$ws.Range("P72").Value = 'Hdac1_CDS_4'
$ws.Range("P73").Value = 'Suv39h1_3UTR_3'
$ws.Range("P74").Value = 'Suv39h1_CDS_4'
$ws.Range("P75").Value = 'Zfp57_CDS_4'
$ws.Range("P76").Value = 'Smchd1_3UTR_5'
$ws.Range("P77").Value = 'Atf7ip_CDS_4'
$ws.Range("P78").Value = 'Ring1_CDS_4'
$ws.Range("P79").Value = 'Atf7ip_CDS_5'
$ws.Range("P80").Value = 'Eed_CDS_5'
$ws.Range("P81").Value = 'Ctcf_CDS_5'
$ws.Range("P82").Value = 'Sirt1_CDS_5'
$ws.Range("P83").Value = 'Zfp57_CDS_5'
$ws.Range("P84").Value = 'Hdac2_CDS_5'
$ws.Range("P85").Value = 'Kmt5a_3UTR_5'
$ws.Range("P86").Value = 'Ring1_CDS_5'
$ws.Range("P87").Value = 'Atrx_CDS_4'
$ws.Range("P88").Value = 'Parp1_CDS_5'
$ws.Range("P89").Value = 'Rnf2_CDS_4'
$ws.Range("P90").Value = 'Atrx_CDS_5'
$ws.Range("P91").Value = 'Suv39h1_CDS_5'
$ws.Range("P92").Value = 'Hdac1_3UTR_5'
$ws.Range("P93").Value = 'Rnf2_CDS_5'
$ws.Range("P94").Value = 'NA_NA_1'
$ws.Range("P95").Value = 'NA_NA_2'
$ws.Range("P96").Value = 'NA_NA_3'
$ws.Range("P97").Value = 'NA_NA_4'

Write-Host "Updated Construct ID (column P) values for rows 1-97"
